$d = $word.ActiveDocument

# Remove the two "Author" style paragraphs ("Laura Kennedy" and "Ben Jarman")
# that directly follow the title. Iterate from the end so removing paragraphs
# doesn't disturb the indices of ones not yet processed.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $styleName = $para.Style.NameLocal
    if ($styleName -eq "Author") {
        $para.Range.Delete()
    }
}
